$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab) from SCD0181 to SCD0011
$ws.Name = "SCD0011"

# Update the TC_ID cell (B2) from "DGS-196" to "SCD0011-012"
$ws.Range("B2").Value = "SCD0011-012"

# Move the active cell selection to B3 (no scrolled top-left cell override)
$ws.Range("B3").Select()

# Column B auto-fits wider to accommodate the new, longer TC_ID text
$ws.Columns.Item(2).ColumnWidth = 11.6
